$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New scrape timestamp applied to every data row.
# ---------------------------------------------------------------------
$ts = "2025-12-15 06:32:13"

# ---------------------------------------------------------------------
# Drop any existing hyperlinks up front; they'll be rebuilt (in row
# order) once all the row data below is in place, so relationship ids
# line up with the final F2..F11 layout.
# ---------------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete()

# ---------------------------------------------------------------------
# Row 2: same listing as before, just refreshed timestamp.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = $ts
$ws.Range("B2").Value = "【AIシステム構築】次のテストに向けた宿題自動出題システム"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5453785"
$ws.Range("G2").Value = 318
$ws.Range("H2").Value = "🔥AI,Ai"

# ---------------------------------------------------------------------
# Row 3: new listing (センサー画像解析 AI案件).
# ---------------------------------------------------------------------
$ws.Range("A3").Value = $ts
$ws.Range("B3").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G3").Value = 310
$ws.Range("H3").Value = "🔥AI,Ai"

# ---------------------------------------------------------------------
# Row 4: new listing (AIオートメーションエンジニア).
# ---------------------------------------------------------------------
$ws.Range("A4").Value = $ts
$ws.Range("B4").Value = "AIオートメーションエンジニア"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5453810"
$ws.Range("G4").Value = 298
$ws.Range("H4").Value = "🔥AI,Ai"

# ---------------------------------------------------------------------
# Row 5: new listing (Unity/XRエンジニア募集).
# ---------------------------------------------------------------------
$ws.Range("A5").Value = $ts
$ws.Range("B5").Value = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5454210"
$ws.Range("G5").Value = 108
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# ---------------------------------------------------------------------
# Row 6: new listing (Javaプログラミング研修).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = $ts
$ws.Range("B6").Value = "Javaプログラミング研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5453723"
$ws.Range("G6").Value = 85
$ws.Range("H6").Value = "★Java"

# ---------------------------------------------------------------------
# Row 7: new listing (GoogleCloud SREエンジニア).
# ---------------------------------------------------------------------
$ws.Range("A7").Value = $ts
$ws.Range("B7").Value = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5453768"
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = "◆開発"

# ---------------------------------------------------------------------
# Row 8: previously row 3 (Base無在庫ツール), pushed down.
# ---------------------------------------------------------------------
$ws.Range("A8").Value = $ts
$ws.Range("B8").Value = "Base無在庫ツール作成 経験者のみ募集 実績提示をお願いします"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5453611"
$ws.Range("G8").Value = 73
$ws.Range("H8").Value = "◆ツール"

# ---------------------------------------------------------------------
# Row 9: new listing (クラウド運用管理研修).
# ---------------------------------------------------------------------
$ws.Range("A9").Value = $ts
$ws.Range("B9").Value = "クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】(再掲)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5453718"
$ws.Range("G9").Value = 38
$ws.Range("H9").Value = "◇管理"

# ---------------------------------------------------------------------
# Row 10: new listing (グーグルワークスペース構築). No H value.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = $ts
$ws.Range("B10").Value = "【介護事業所向け】グーグルワークスペース社内システム構築依頼"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5453868"
$ws.Range("G10").Value = 40
$ws.Range("H10").ClearContents()

# ---------------------------------------------------------------------
# Row 11: previously row 4 (ホームページとLP), pushed down. No H value.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = $ts
$ws.Range("B11").Value = "【急募】ホームページとLPの改善をお手伝いします!"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5453763"
$ws.Range("G11").Value = 18
$ws.Range("H11").ClearContents()

# ---------------------------------------------------------------------
# Rebuild hyperlinks on the URL column, in row order, so relationship
# ids come out rId1..rId10 matching F2..F11.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5453785")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5427956")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5453810")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5454210")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5453723")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5453768")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5453611")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5453718")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5453868")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5453763")

# ---------------------------------------------------------------------
# Column width tweaks (B: 34 -> 51, D: 28 -> 30). ColumnWidth persists
# with Excel's +5/6-character padding baked into the stored value, so
# back that off here to land on the exact target width.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 51 - 5/6
$ws.Columns.Item(4).ColumnWidth = 30 - 5/6
